$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

# Set "Create Test Passed" (column B) to TRUE for rows 2-15
$ws.Range("B2:B15").Value = $true
